# "complete alarm for sim"
# Adds two new rows (hkck_11 / hkck_33 - Hospital Sg Buloh, a Covid-era
# hospital critical-list entry that was removed in 2022) to the
# flaglist_incomer sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the 33kV trip id first so it lands in the shared-string table
# ahead of the other new strings (matches author's original edit order).
$ws.Cells.Item(65, 1).Value = "hkck_33"
$ws.Cells.Item(65, 3).Value = "Hospital Sg Buloh"

$ws.Cells.Item(64, 1).Value = "hkck_11"
$ws.Cells.Item(64, 4).Value = "Covid Hosipital - removed on 2022 during covid pandemik"

$ws.Cells.Item(64, 2).Value = "gso"
$ws.Cells.Item(64, 3).Value = "Hospital Sg Buloh"

$ws.Cells.Item(65, 2).Value = "gso"
$ws.Cells.Item(65, 4).Value = "Covid Hosipital - removed on 2022 during covid pandemik"

# Column D has no sheet-level default style, so the long_text cells in
# rows 64/65 need their formatting copied explicitly from the row above
# (style index 15 - Aptos Narrow 10pt) to match the rest of the column.
$ws.Range("D63").Copy()
$ws.Range("D64").PasteSpecial(-4122)
$ws.Range("D65").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$null = $ws.Range("H69").Select()
